$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Helper: set a text value on a cell while preventing Excel's
# automatic "looks like a date" -> date-serial conversion. We do
# this by temporarily copying the format from a reference cell,
# switching to a text number format, assigning the value, then
# restoring the original format by pasting the reference cell's
# format back on top (keeps the exact same style / numFmtId).
# ---------------------------------------------------------------
function Set-TextValue($cellRef, $formatSourceRef, $value) {
    $ws.Range($formatSourceRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
    $ws.Range($formatSourceRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

# -----------------------------------------------------------------
# Update existing rows 2 & 3 - new payment contact name/email values
# -----------------------------------------------------------------
$ws.Range("R2").Value = "Monika Single Individual"
$ws.Range("T2").Value = "monika.srivastava+mssingleindividualauto1@aidenai.com"

Set-TextValue "J3" "J3" "05-20-2024"
$ws.Range("R3").Value = "Monika Multiple Individual"
$ws.Range("T3").Value = "monika.srivastava+msmultipleindividualauto1@aidenai.com"

# Re-point existing hyperlinks' targets to the updated e-mail addresses
$ws.Hyperlinks.Item(2).Address = "mailto:monika.srivastava+mssingleindividualauto1@aidenai.com"
$ws.Hyperlinks.Item(4).Address = "mailto:monika.srivastava+msmultipleindividualauto1@aidenai.com"

# -----------------------------------------------------------------
# Row 4 - new scenario: validate_Event_E2EScenario_SingleEvent_Business
# -----------------------------------------------------------------
$ws.Range("A4").Value = "validate_Event_E2EScenario_SingleEvent_Business"
$ws.Range("B4").Value = "Bridal Shower"
$ws.Range("C4").Value = "403 Church Street, Toronto, ON, Canada"
$ws.Range("D4").Value = "Single event"
$ws.Range("E4").Value = 5
Set-TextValue "F4" "F3" "05-20-2024"
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = ":15"
$ws.Range("I4").Value = "AM"
Set-TextValue "J4" "J2" "05-25-2024"
$ws.Range("K4").Value = 11
$ws.Range("L4").Value = ":30"
$ws.Range("M4").Value = "PM"
$ws.Range("N4").Value = "string:between 26-125"
$ws.Range("O4").Value = "No"
$ws.Range("P4").Value = "Bridal Shower"
$ws.Range("Q4").Value = "A Business"
$ws.Range("R4").Value = "Monika Single Business"
$ws.Range("S4").Value = "string:Ontario"
$ws.Range("T4").Value = "monika.srivastava+mssinglebusinessauto1@aidenai.com"
$ws.Range("U4").Value = 1557773334
$ws.Range("V4").Value = "Event_Single_UI_Business"

# -----------------------------------------------------------------
# Row 5 - new scenario: validate_Event_E2EScenario_MultipleEvents_Business
# -----------------------------------------------------------------
$ws.Range("A5").Value = "validate_Event_E2EScenario_MultipleEvents_Business"
$ws.Range("B5").Value = "Auction"
$ws.Range("C5").Value = "3032 Dougall Avenue, Windsor, ON, Canada"
$ws.Range("D5").Value = "Multiple events"
$ws.Range("E5").Value = 10
Set-TextValue "F5" "F3" "05-30-2024"
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = ":30"
$ws.Range("I5").Value = "PM"
Set-TextValue "J5" "J2" "06-30-2024"
$ws.Range("K5").Value = 9
$ws.Range("L5").Value = ":40"
$ws.Range("M5").Value = "PM"
$ws.Range("N5").Value = "string:between 501-1,000"
$ws.Range("O5").Value = "No"
$ws.Range("P5").Value = "Auction "
$ws.Range("Q5").Value = "A Business"
$ws.Range("R5").Value = "Monika Multiple Business"
$ws.Range("S5").Value = "string:Ontario"
$ws.Range("T5").Value = "monika.srivastava+msmultiplebusinessauto1@aidenai.com"
$ws.Range("U5").Value = 1555999333
$ws.Range("V5").Value = "Event_Multiple_UI_Business"

# -----------------------------------------------------------------
# Hyperlinks for the two new rows (mirrors the rows 2/3 pattern:
# "O" column links to the fixed single-individual mailbox, "T"
# column links to the row's own entity e-mail address)
# -----------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("O4"), "mailto:monika.srivastava+eventsingleindividual@aidenai.com", "", "", "monika.srivastava+eventsingleindividual@aidenai.com") | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null

$ws.Hyperlinks.Add($ws.Range("T4"), "mailto:monika.srivastava+mssinglebusinessauto1@aidenai.com", "", "", "monika.srivastava+mssinglebusinessauto1@aidenai.com") | Out-Null
$ws.Range("T3").Copy() | Out-Null
$ws.Range("T4").PasteSpecial(-4122) | Out-Null
$ws.Range("T4").Value = "monika.srivastava+mssinglebusinessauto1@aidenai.com"

$ws.Hyperlinks.Add($ws.Range("O5"), "mailto:monika.srivastava+eventsingleindividual@aidenai.com", "", "", "monika.srivastava+eventsingleindividual@aidenai.com") | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null

$ws.Hyperlinks.Add($ws.Range("T5"), "mailto:monika.srivastava+msmultiplebusinessauto1@aidenai.com", "", "", "monika.srivastava+msmultiplebusinessauto1@aidenai.com") | Out-Null
$ws.Range("T3").Copy() | Out-Null
$ws.Range("T5").PasteSpecial(-4122) | Out-Null
$ws.Range("T5").Value = "monika.srivastava+msmultiplebusinessauto1@aidenai.com"

# -----------------------------------------------------------------
# Column V got a touch wider (Excel best-fit) after the new content
# was added.
# -----------------------------------------------------------------
$ws.Range("V1").ColumnWidth = 24.3

# -----------------------------------------------------------------
# Final active selection, as recorded in the saved workbook.
# -----------------------------------------------------------------
$ws.Range("J11").Select() | Out-Null
